$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "64.243.13"
Set-TextCell 2 5 "  -3.25%  "

# Row 3
Set-TextCell 3 4 "3.167.17"
Set-TextCell 3 5 "  -4.55%  "

# Row 4
Set-TextCell 4 5 "  +0.03%  "

# Row 5
Set-TextCell 5 4 "568.53"
Set-TextCell 5 5 "  -3.07%  "

# Row 6
Set-TextCell 6 4 "169.03"
Set-TextCell 6 5 "  -7.45%  "

# Row 7
Set-TextCell 7 4 "0.609"
Set-TextCell 7 5 "  -5.49%  "

# Row 8
Set-TextCell 8 5 "  +0.06%  "

# Row 9
Set-TextCell 9 4 "3.167.80"
Set-TextCell 9 5 "  -4.44%  "

# Row 10
Set-TextCell 10 4 "0.121"
Set-TextCell 10 5 "  -4.48%  "

# Row 11
Set-TextCell 11 4 "6.72"
Set-TextCell 11 5 "  -1.00%  "

# Row 12
Set-TextCell 12 5 "  -4.57%  "

# Row 13
Set-TextCell 13 4 "3.724.46"
Set-TextCell 13 5 "  -4.35%  "

# Row 14
Set-TextCell 14 5 "  -2.02%  "

# Row 15
Set-TextCell 15 4 "64.281.92"
Set-TextCell 15 5 "  -3.20%  "

# Row 16
Set-TextCell 16 4 "25.43"
Set-TextCell 16 5 "  -3.79%  "

# Row 17
Set-TextCell 17 5 "  -3.77%  "

# Row 18
Set-TextCell 18 4 "3.180.49"
Set-TextCell 18 5 "  -3.98%  "

# Row 19
Set-TextCell 19 4 "419.68"
Set-TextCell 19 5 "  -2.59%  "

# Row 20
Set-TextCell 20 2 "Chainlink"
Set-TextCell 20 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 20 4 "12.80"
Set-TextCell 20 5 "  -3.94%  "

# Row 21
Set-TextCell 21 2 "Polkadot"
Set-TextCell 21 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 21 4 "5.34"
Set-TextCell 21 5 "  -3.53%  "

# Row 22
Set-TextCell 22 4 "7.03"
Set-TextCell 22 5 "  -5.32%  "

# Row 23
Set-TextCell 23 5 "  -0.20%  "

# Row 24
Set-TextCell 24 4 "70.20"
Set-TextCell 24 5 "  -2.95%  "

# Row 25
Set-TextCell 25 4 "0.202"
Set-TextCell 25 5 "  +2.53%  "

# Row 26
Set-TextCell 26 4 "0.489"
Set-TextCell 26 5 "  -5.33%  "

# Row 28
Set-TextCell 28 4 "8.74"
Set-TextCell 28 5 "  -3.23%  "

# Row 29
Set-TextCell 29 4 "1.00"
Set-TextCell 29 5 "  +0.01%  "

# Row 30
Set-TextCell 30 4 "21.80"
Set-TextCell 30 5 "  -2.87%  "

# Row 31
Set-TextCell 31 5 "  -6.62%  "

# Row 32
Set-TextCell 32 4 "0.998"
Set-TextCell 32 5 "  -0.12%  "

# Row 33
Set-TextCell 33 5 "  -4.49%  "

# Row 34
Set-TextCell 34 4 "6.32"
Set-TextCell 34 5 "  -4.80%  "

# Row 35
Set-TextCell 35 5 "  -5.96%  "

# Row 36
Set-TextCell 36 4 "156.73"
Set-TextCell 36 5 "  -1.45%  "

# Row 37
Set-TextCell 37 5 "  -7.03%  "

# Row 38
Set-TextCell 38 4 "2.703.34"
Set-TextCell 38 5 "  -6.01%  "

# Row 39
Set-TextCell 39 5 "  -7.50%  "

# Row 40
Set-TextCell 40 4 "24.46"
Set-TextCell 40 5 "  -8.91%  "

# Row 41
Set-TextCell 41 5 "  -4.22%  "

# Row 42
Set-TextCell 42 4 "39.07"
Set-TextCell 42 5 "  -2.93%  "

# Row 43
Set-TextCell 43 4 "0.707"
Set-TextCell 43 5 "  -8.23%  "

# Row 44
Set-TextCell 44 4 "5.72"
Set-TextCell 44 5 "  -5.40%  "

# Row 45
Set-TextCell 45 4 "0.0621"
Set-TextCell 45 5 "  -6.88%  "

# Row 46
Set-TextCell 46 2 "VeChain"
Set-TextCell 46 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 46 4 "0.0261"
Set-TextCell 46 5 "  -3.71%  "

# Row 47
Set-TextCell 47 2 "InjectiveProtocol"
Set-TextCell 47 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell 47 4 "21.74"
Set-TextCell 47 5 "  -7.34%  "

# Row 48
Set-TextCell 48 2 "Bittensor"
Set-TextCell 48 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell 48 4 "294.17"
Set-TextCell 48 5 "  -7.48%  "

# Row 49
Set-TextCell 49 2 "FirstDigitalUSD"
Set-TextCell 49 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell 49 4 "1.00"
Set-TextCell 49 5 "  +0.03%  "

# Row 50
Set-TextCell 50 2 "dogwifhat"
Set-TextCell 50 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell 50 4 "2.01"
Set-TextCell 50 5 "  -13.43%  "

# Row 51
Set-TextCell 51 4 "0.0992"
Set-TextCell 51 5 "  -4.46%  "
